$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price values (column D).
# NumberFormat is forced to text ("@") per cell before the write so
# values such as "1.00" / "212.11" persist as the literal text that
# the source site renders, instead of being reinterpreted as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.662.69"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.603.77"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.11"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.516"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.02"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.42"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.833.04"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.590.35"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.660.98"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.36"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.76"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.57"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.08"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.72"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.427.95"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.537"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "54.95"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.95"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.89"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.954"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.742.57"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "86.46"

# Updated 1h volume/change percentages (column E)
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +8.45%  "
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("E15").Value = "  +3.43%  "
$ws.Range("E16").Value = "  +3.63%  "
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("E19").Value = "  +5.88%  "
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("E27").Value = "  +3.56%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +3.70%  "
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("E34").Value = "  +4.27%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("E43").Value = "  +1.53%  "
$ws.Range("E44").Value = "  +6.65%  "
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("E47").Value = "  +3.21%  "
$ws.Range("E48").Value = "  +16.00%  "
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("E51").Value = "  +2.09%  "
